# Update column F ("dSF") values on Sheet1 to reflect the repulled/recalculated data.
# Mapping of row -> new value for column F.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    2  = 1
    7  = 1
    11 = -1
    12 = 4
    16 = 0
    24 = 4
    25 = -1
    26 = 1
    30 = 3
    31 = 0
    43 = -1
    44 = 0
    45 = 2
    46 = 0
    47 = 3
    49 = 0
    56 = 1
    71 = -2
    74 = -2
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}
